$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.281.41"
$ws.Range("D3").Value = "1.868.17"
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'235.14"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4697"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "'0.2872"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "'0.06570"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'21.80"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").Value = "'0.08015"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'97.18"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.872.28"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "'5.123"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'0.6850"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "'269.41"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").Value = "30.269.74"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'14.04"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").Value = "'0.000007662"
$ws.Range("E19").Value = "  +5.03%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "2.116.15"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'5.275"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").Value = "'6.209"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "'9.444"
$ws.Range("E25").Value = "  +3.06%  "
$ws.Range("D26").Value = "'168.27"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'18.90"
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "'1.947"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").Value = "'0.09873"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").Value = "'4.375"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "'1.464"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "'4.074"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "'0.7001"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "'2.712"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'0.01874"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'2.626"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'6.294"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'72.33"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").Value = "'1.951"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "'0.8418"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4162"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "'103.02"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.190"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.053"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "'923.12"
$ws.Range("E49").Value = "  -5.97%  "
$ws.Range("D50").Value = "'34.47"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  +0.66%  "
